# Apply crypto price/volume/name/link updates (GitHub Actions refresh).
# Values that Excel would otherwise auto-convert to numbers (e.g. "1.001",
# "0.7331") are written with a leading apostrophe so they stay literal text,
# matching the workbook's inlineStr cell storage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.905.75"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "1.887.12"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'0.7331"
$ws.Range("E5").Value = "  -4.61%  "
$ws.Range("D6").Value = "'242.46"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.3106"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").Value = "'26.22"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").Value = "'0.06895"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.7689"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07940"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.877.37"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "'5.218"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "'91.20"
$ws.Range("E15").Value = "  -3.55%  "
$ws.Range("D16").Value = "'14.17"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "29.941.55"
$ws.Range("D18").Value = "'5.761"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'239.54"
$ws.Range("D20").Value = "'0.000007751"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "2.158.12"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'6.934"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").Value = "'9.275"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "'164.41"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'18.84"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'0.1270"
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("D29").Value = "'2.013"
$ws.Range("E29").Value = "  -11.29%  "
$ws.Range("D30").Value = "'1.359"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "'1.529"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'4.297"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").Value = "'4.082"
$ws.Range("D34").Value = "'0.05085"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "'1.275"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'0.7361"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'2.720"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'0.01918"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'2.773"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'6.298"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").Value = "'74.53"
$ws.Range("E41").Value = "  -5.38%  "
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "'1.930"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "'0.8364"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'7.609"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").Value = "'100.84"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'9.785"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.056.15"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.97"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "'942.88"
$ws.Range("E51").Value = "  -3.55%  "
